# Tutor Profile Student Message
# Applies the edits to the "Notes" workbook:
#  - Sheet "Notes": rewrite rows 5-7, add new rows 8-12 (tutor54 sessions + a
#    new student-message row), refresh related hyperlinks.
#  - Sheet "Notes_list": remove a blank spacer row (17 -> 16) and update the
#    selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Notes"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Notes")

# Make room for 5 new rows (8..12) by copying row 5's formatting (style
# pattern: C=left-aligned, D=hyperlink+fill/border, E=hyperlink, F=plain)
# down into the new rows.
$ws.Rows.Item(5).Copy()
$ws.Range("A8:A12").EntireRow.Insert(-4121)

# F7 loses its italic "Arial 7" style (s=4) in the new layout; new rows
# 8-11 should not carry the style copied from row 5's F (s=4 comes from the
# inserted rows picking up row 7's neighbourhood) - reset every F cell we
# touch to the default (Normal) style, and D12 also reverts to the default
# (unstyled) look.
$ws.Range("F7:F12").Style = "Normal"
$ws.Range("D12").Style = "Normal"

# --- Row 5 -----------------------------------------------------------------
$ws.Range("B5").Value = "d"
$ws.Range("D5").Value = "tutor53@nkt.com"
$ws.Range("F5").Value = "Java session 8"

# --- Row 6 -----------------------------------------------------------------
$ws.Range("D6").Value = "tutor53@nkt.com"
$ws.Range("F6").Value = "Mv class 2"

# --- Row 7 (fully rewritten) -------------------------------------------------
$ws.Range("A7").Value = "signin"
$ws.Range("B7").Value = "d"
$ws.Range("C7").Value = "n"
$ws.Range("D7").Value = "tutor54@nkt.com"
$ws.Range("E7").Value = "Admin@123"
$ws.Range("F7").Value = "paint sess 1"
$ws.Range("G7").Value = 4

# --- Row 8 (new) -------------------------------------------------------------
$ws.Range("A8").Value = "signin"
$ws.Range("B8").Value = "d"
$ws.Range("C8").Value = "n"
$ws.Range("D8").Value = "tutor54@nkt.com"
$ws.Range("E8").Value = "Admin@123"
$ws.Range("F8").Value = "paint fix 1"
$ws.Range("G8").Value = 4

# --- Row 9 (new) -------------------------------------------------------------
$ws.Range("A9").Value = "signin"
$ws.Range("B9").Value = "d"
$ws.Range("C9").Value = "n"
$ws.Range("D9").Value = "tutor54@nkt.com"
$ws.Range("E9").Value = "Admin@123"
$ws.Range("F9").Value = "paint var 1"
$ws.Range("G9").Value = 4

# --- Row 10 (new) ------------------------------------------------------------
$ws.Range("A10").Value = "signin"
$ws.Range("B10").Value = "tutor"
$ws.Range("C10").Value = "n"
$ws.Range("D10").Value = "tutor54@nkt.com"
$ws.Range("E10").Value = "Admin@123"
$ws.Range("F10").Value = "paint sess 1 ind"
$ws.Range("G10").Value = 4

# --- Row 11 (new) ------------------------------------------------------------
$ws.Range("A11").Value = "signin"
$ws.Range("B11").Value = "tutor"
$ws.Range("C11").Value = "n"
$ws.Range("D11").Value = "tutor54@nkt.com"
$ws.Range("E11").Value = "Admin@123"
$ws.Range("F11").Value = "paint fix 1 ind"
$ws.Range("G11").Value = 4

# --- Row 12 (new - student message test case) --------------------------------
$ws.Range("A12").Value = "signin"
$ws.Range("B12").Value = "d"
$ws.Range("C12").Value = "n"
$ws.Range("D12").Value = "srinivasesaivanan6324@gmail.com"
$ws.Range("E12").Value = "Test@1234"
$ws.Range("F12").Value = "clarinet session 4 multi"
$ws.Range("G12").Value = 4

# --- Hyperlinks: rebuild the full list so ids/order match the new layout ----
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:tutor39@nkt.com")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:tutor53@nkt.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:tutor53@nkt.com")
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:tutor54@nkt.com")
$ws.Hyperlinks.Add($ws.Range("E8"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("E9"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("E10"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("E11"), "mailto:Admin@123")
$ws.Hyperlinks.Add($ws.Range("E12"), "mailto:Test@1234")

$ws.Range("B12").Select()

# ---------------------------------------------------------------------------
# Sheet "Notes_list"
# ---------------------------------------------------------------------------
$wsList = $wb.Worksheets.Item("Notes_list")

# Remove the blank spacer row (was row 4), shifting the trailing note
# (previously row 17) up to row 16.
$wsList.Rows.Item(4).Delete()

$wsList.Rows.Item(4).Select()
